$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh coin price / 1h change figures from the latest scrape.
# Rows 37-38 also swap rank: Monero moved above ImmutableX.
#
# Some Price values are plain decimals (e.g. "583.15") that Excel
# would otherwise auto-convert to a real number. A leading apostrophe
# forces them to stay text, like the original cells; Style is then
# reset to "Normal" to drop the quote-prefix mark Excel adds so the
# cell keeps its original (default) formatting.

$ws.Range("D2").Value = '66.379.51'
$ws.Range("E2").Value = '  -1.15%  '
$ws.Range("D3").Value = '2.559.85'
$ws.Range("E3").Value = '  -2.33%  '
$ws.Range("E4").Value = '  +0.14%  '
$r = $ws.Range("D5")
$r.Value = "'583.15"
$r.Style = "Normal"
$ws.Range("E5").Value = '  -1.90%  '
$r = $ws.Range("D6")
$r.Value = "'167.26"
$r.Style = "Normal"
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("E7").Value = '  +0.22%  '
$r = $ws.Range("D8")
$r.Value = "'0.526"
$r.Style = "Normal"
$ws.Range("E8").Value = '  -1.20%  '
$ws.Range("D9").Value = '2.558.62'
$ws.Range("E9").Value = '  -2.39%  '
$r = $ws.Range("D10")
$r.Value = "'0.139"
$r.Style = "Normal"
$ws.Range("E10").Value = '  -0.62%  '
$ws.Range("E11").Value = '  +0.19%  '
$r = $ws.Range("D12")
$r.Value = "'0.354"
$r.Style = "Normal"
$ws.Range("E12").Value = '  -2.03%  '
$r = $ws.Range("D13")
$r.Value = "'5.14"
$r.Style = "Normal"
$ws.Range("E13").Value = '  -1.68%  '
$r = $ws.Range("D14")
$r.Value = "'26.64"
$r.Style = "Normal"
$ws.Range("E14").Value = '  -3.47%  '
$ws.Range("D15").Value = '3.029.25'
$ws.Range("E15").Value = '  -2.52%  '
$ws.Range("E16").Value = '  -2.23%  '
$ws.Range("D17").Value = '66.247.01'
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '2.541.93'
$ws.Range("E18").Value = '  -3.01%  '
$r = $ws.Range("D19")
$r.Value = "'11.37"
$r.Style = "Normal"
$ws.Range("E19").Value = '  -6.23%  '
$r = $ws.Range("D20")
$r.Value = "'7.73"
$r.Style = "Normal"
$ws.Range("E20").Value = '  -4.32%  '
$r = $ws.Range("D21")
$r.Value = "'349.43"
$r.Style = "Normal"
$ws.Range("E21").Value = '  -1.96%  '
$r = $ws.Range("D22")
$r.Value = "'4.21"
$r.Style = "Normal"
$ws.Range("E22").Value = '  -2.55%  '
$r = $ws.Range("D23")
$r.Value = "'4.58"
$r.Style = "Normal"
$ws.Range("E23").Value = '  -1.77%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  +0.19%  '
$r = $ws.Range("D26")
$r.Value = "'69.20"
$r.Style = "Normal"
$ws.Range("E26").Value = '  -1.21%  '
$r = $ws.Range("D27")
$r.Value = "'9.90"
$r.Style = "Normal"
$ws.Range("E27").Value = '  -6.61%  '
$ws.Range("D28").Value = '2.696.42'
$ws.Range("E28").Value = '  -2.24%  '
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("D30").Value = '0.0₃0989'
$ws.Range("E30").Value = '  -1.43%  '
$r = $ws.Range("D31")
$r.Value = "'8.27"
$r.Style = "Normal"
$ws.Range("E31").Value = '  +4.29%  '
$r = $ws.Range("D32")
$r.Value = "'526.27"
$r.Style = "Normal"
$ws.Range("E32").Value = '  -4.03%  '
$ws.Range("E33").Value = '  -2.51%  '
$ws.Range("E34").Value = '  -3.94%  '
$r = $ws.Range("D35")
$r.Value = "'0.131"
$r.Style = "Normal"
$ws.Range("E35").Value = '  -3.53%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$r = $ws.Range("D37")
$r.Value = "'157.29"
$r.Style = "Normal"
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$r = $ws.Range("D38")
$r.Value = "'1.45"
$r.Style = "Normal"
$ws.Range("E38").Value = '  -3.49%  '
$r = $ws.Range("D39")
$r.Value = "'18.71"
$r.Style = "Normal"
$ws.Range("E39").Value = '  -1.54%  '
$r = $ws.Range("D40")
$r.Value = "'0.357"
$r.Style = "Normal"
$ws.Range("E40").Value = '  -2.29%  '
$r = $ws.Range("D41")
$r.Value = "'18.31"
$r.Style = "Normal"
$ws.Range("E41").Value = '  +2.31%  '
$ws.Range("E42").Value = '  -1.32%  '
$r = $ws.Range("D43")
$r.Value = "'5.08"
$r.Style = "Normal"
$ws.Range("E43").Value = '  -2.28%  '
$ws.Range("E44").Value = '  +0.07%  '
$r = $ws.Range("D45")
$r.Value = "'2.42"
$r.Style = "Normal"
$ws.Range("E45").Value = '  +0.84%  '
$ws.Range("D46").Value = '0.0₆0285'
$ws.Range("E46").Value = '  -3.80%  '
$r = $ws.Range("D47")
$r.Value = "'148.18"
$r.Style = "Normal"
$ws.Range("E47").Value = '  -2.09%  '
$r = $ws.Range("D48")
$r.Value = "'0.563"
$r.Style = "Normal"
$ws.Range("E48").Value = '  -2.45%  '
$ws.Range("E49").Value = '  -2.01%  '
$ws.Range("E50").Value = '  +1.45%  '
$r = $ws.Range("D51")
$r.Value = "'0.0759"
$r.Style = "Normal"
$ws.Range("E51").Value = '  -1.36%  '
